$wb = $excel.ActiveWorkbook

# Weekly update: add data for 2023-11-05 (full-year 2023 running totals, column J)

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6530
$ws.Range("J3").Value = 6924
$ws.Range("C4").Value = 1837
$ws.Range("D4").Value = 1960
$ws.Range("H4").Value = 1707
$ws.Range("I4").Value = 1775
$ws.Range("J4").Value = 1503
$ws.Range("J5").Value = 535
$ws.Range("J6").Value = 9189
$ws.Range("C7").Value = 28381
$ws.Range("D7").Value = 28150
$ws.Range("H7").Value = 26018
$ws.Range("I7").Value = 26232
$ws.Range("J7").Value = 24681

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 197
$ws.Range("J6").Value = 188
$ws.Range("J7").Value = 720
$ws.Range("J8").Value = 1551
$ws.Range("J9").Value = 132
$ws.Range("J10").Value = 180
$ws.Range("J11").Value = 421
$ws.Range("J12").Value = 52
$ws.Range("J13").Value = 30
$ws.Range("J15").Value = 292
$ws.Range("J19").Value = 720
$ws.Range("J20").Value = 518
$ws.Range("J23").Value = 226
$ws.Range("J25").Value = 121
$ws.Range("J28").Value = 7
$ws.Range("J29").Value = 1343
$ws.Range("J31").Value = 241
$ws.Range("J33").Value = 1112
$ws.Range("J37").Value = 754
$ws.Range("J41").Value = 171
$ws.Range("J42").Value = 1061
$ws.Range("J43").Value = 212
$ws.Range("J44").Value = 187
$ws.Range("J46").Value = 82
$ws.Range("J51").Value = 300
$ws.Range("J52").Value = 621
$ws.Range("J53").Value = 356
$ws.Range("J54").Value = 467
$ws.Range("J55").Value = 375
$ws.Range("J60").Value = 144
$ws.Range("C63").Value = 267
$ws.Range("D63").Value = 344
$ws.Range("H63").Value = 265
$ws.Range("I63").Value = 252
$ws.Range("J63").Value = 88
$ws.Range("J65").Value = 612
$ws.Range("J67").Value = 935
$ws.Range("J71").Value = 80
$ws.Range("J75").Value = 74
$ws.Range("J78").Value = 293
$ws.Range("J79").Value = 693
$ws.Range("J83").Value = 489
$ws.Range("J85").Value = 1023
$ws.Range("J87").Value = 81
$ws.Range("J88").Value = 253
$ws.Range("J90").Value = 263
$ws.Range("J96").Value = 271
$ws.Range("J97").Value = 220
$ws.Range("J98").Value = 182
$ws.Range("C101").Value = 28381
$ws.Range("D101").Value = 28150
$ws.Range("H101").Value = 26018
$ws.Range("I101").Value = 26232
$ws.Range("J101").Value = 24681

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 81
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 218
$ws.Range("J6").Value = 234
$ws.Range("J7").Value = 720

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 121
$ws.Range("J3").Value = 79
$ws.Range("J6").Value = 190
$ws.Range("J7").Value = 421

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 361
$ws.Range("J7").Value = 1023

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 143
$ws.Range("J7").Value = 621

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 236
$ws.Range("J7").Value = 356

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 423
$ws.Range("J3").Value = 465
$ws.Range("J6").Value = 542
$ws.Range("J7").Value = 1551

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 183
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 489

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 368
$ws.Range("J4").Value = 48
$ws.Range("J5").Value = 46
$ws.Range("J6").Value = 389
$ws.Range("J7").Value = 1112

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 221
$ws.Range("J6").Value = 223
$ws.Range("J7").Value = 754

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J5").Value = 17
$ws.Range("J7").Value = 612

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 86
$ws.Range("J3").Value = 63
$ws.Range("J4").Value = 12
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 241

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 237
$ws.Range("J3").Value = 347
$ws.Range("J6").Value = 260
$ws.Range("J7").Value = 935

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 97
$ws.Range("J6").Value = 219
$ws.Range("J7").Value = 467

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 409
$ws.Range("J3").Value = 472
$ws.Range("J6").Value = 339
$ws.Range("J7").Value = 1343

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 170
$ws.Range("J7").Value = 720

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 75
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 57
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 25
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 225
$ws.Range("J3").Value = 207
$ws.Range("J6").Value = 566
$ws.Range("J7").Value = 1061

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J3").Value = 9
$ws.Range("J6").Value = 30

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 293

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 210
$ws.Range("J7").Value = 375

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 60
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 235
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 693

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 146
$ws.Range("J4").Value = 42
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 518

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 66
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 125
$ws.Range("J7").Value = 292

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 114
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J4").Value = 11
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 119
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 7
